# Add files via upload
# Extend the "Source" sheet data from columns B:I out to B:N (duplicating
# the last 5 data columns with some updated figures) for rows 2-4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 ("sprint")
$ws.Range("B2").Value = 30
$ws.Range("C2").Value = 45
$ws.Range("D2").Value = 3
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = 5
$ws.Range("G2").Value = 6
$ws.Range("H2").Value = 7
$ws.Range("I2").Value = 10
$ws.Range("J2").Value = 4
$ws.Range("K2").Value = 5
$ws.Range("L2").Value = 6
$ws.Range("M2").Value = 7
$ws.Range("N2").Value = 10

# Row 3 ("Duration")
$ws.Range("B3").Value = 60
$ws.Range("C3").Value = 28
$ws.Range("D3").Value = 15
$ws.Range("E3").Value = 7
$ws.Range("F3").Value = 17
$ws.Range("G3").Value = 100
$ws.Range("H3").Value = 15
$ws.Range("I3").Value = 20
$ws.Range("J3").Value = 7
$ws.Range("K3").Value = 17
$ws.Range("L3").Value = 100
$ws.Range("M3").Value = 15
$ws.Range("N3").Value = 20

# Row 4 ("Grooming")
$ws.Range("B4").Value = 4
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 90
$ws.Range("E4").Value = 150
$ws.Range("F4").Value = 3
$ws.Range("G4").Value = 37
$ws.Range("H4").Value = 0.5
$ws.Range("I4").Value = 130
$ws.Range("J4").Value = 150
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 37
$ws.Range("M4").Value = 0.5
$ws.Range("N4").Value = 130

# Apply the same header fill/format (style 2, from B1) to the new header
# cells J1:N1 so the look is consistent with the rest of row 1.
$ws.Range("J1:N1").Value = $ws.Range("I1").Value

# Match the new cell formatting (style 3, center aligned w/ fill) used by
# the rest of the data block for the newly added columns.
$ws.Range("J2:N4").Style = $ws.Range("I2").Style

# Reflect the new used range / selection as captured by the workbook.
$ws.Range("N2").Select
$excel.ActiveWindow.ScrollColumn = 2
